$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update H34 and H38 from 0 to -1
$ws.Range("H34").Value = -1
$ws.Range("H38").Value = -1

# Delete row 48 (Maja squinado), shifting rows 49-65 up
$ws.Rows("48").Delete()
